$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before the current row 10 ("Het is mogelijk...") for the
# reworded "Er kan een wapen gekozen worden..." requirement (now with the
# "(waarden tussen 0 en 100)" qualifier added after review).
$ws.Rows.Item(10).Insert()
$ws.Cells.Item(10, 1).Value = "Er kan een wapen gekozen worden op basis van schade per kogel (waarden tussen 0 en 100)"
$ws.Cells.Item(10, 2).Value = "Must have"

# Insert a new row after that for the new "Counter die de speltijd..." requirement.
$ws.Rows.Item(12).Insert()
$ws.Cells.Item(12, 1).Value = "Counter die de speltijd bijhoudt en vergelijkt met de eerder ingestelde totale speltijd, wanneer deze gelijk zijn betekend dit het einde van het spel"
$ws.Cells.Item(12, 2).Value = "Must have"

# The old, un-qualified "Er kan een wapen gekozen worden op basis van schade
# per kogel" row has now shifted down to row 13 (superseded by the reworded
# version above) -- remove it.
$ws.Rows.Item(13).Delete()

# Match the author's final selection.
$ws.Range("A10:B10").Select()
